# Apply the updated crypto price/volume figures captured by the
# scheduled GitHub Actions refresh (commit: "Updated cryptos list").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---
# Several prices are plain decimal-looking strings (e.g. "1.002", "5.890")
# that Excel would otherwise auto-coerce into numbers, dropping the
# significant trailing zeros. Force the cell to Text format first, write
# the literal string, then drop the style back to Normal so no stray
# cell-style index lingers on cells that previously had none.
$ws.Range("D2").Value = "27.452.40"
$ws.Range("D3").Value = "1.750.48"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "322.37"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4252"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3601"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "42.41"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07462"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.68"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.024"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.204"
$c.Style = "Normal"
$ws.Range("D16").Value = "1.746.37"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "93.12"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06369"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.889"
$c.Style = "Normal"
$ws.Range("D23").Value = "27.500.55"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.21"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.086"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "161.99"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.19"
$c.Style = "Normal"
$ws.Range("D28").Value = "1.944.89"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.137"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "123.75"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.103"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.664"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.560"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.08888"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02292"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.2092"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05999"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.957"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "7.924"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.389"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.37"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5875"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.690"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "123.20"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.965"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.162"
$c.Style = "Normal"

# --- Volume(1h) column (E) updates ---
# These already contain "%" and padding spaces, so Excel keeps them as
# plain text without any extra coercion handling.
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -4.53%  "
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("E9").Value = "  -5.84%  "
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("E11").Value = "  -3.47%  "
$ws.Range("E13").Value = "  -6.87%  "
$ws.Range("E14").Value = "  -4.83%  "
$ws.Range("E15").Value = "  -5.69%  "
$ws.Range("E16").Value = "  -5.88%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("E22").Value = "  -6.02%  "
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("E24").Value = "  -4.33%  "
$ws.Range("E25").Value = "  -5.09%  "
$ws.Range("E26").Value = "  +3.39%  "
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("E28").Value = "  -5.19%  "
$ws.Range("E29").Value = "  -8.21%  "
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("E31").Value = "  -9.15%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  -6.75%  "
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("E35").Value = "  -8.06%  "
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("E39").Value = "  -4.68%  "
$ws.Range("E40").Value = "  -5.00%  "
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("E45").Value = "  -4.78%  "
$ws.Range("E46").Value = "  -4.75%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("E49").Value = "  -4.26%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("E51").Value = "  -2.64%  "
